$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.383.69'
$ws.Range("E2").Value = '  +3.35%  '

$ws.Range("D3").Value = '1.868.79'
$ws.Range("E3").Value = '  +1.84%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.23%  '

$ws.Range("E6").Value = '  -0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4698'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3962'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.66%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.63'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08032'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.87%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9995'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.96'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.07%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.872.75'
$ws.Range("E13").Value = '  +1.32%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.031'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.253'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.51%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.08%  '

$ws.Range("E18").Value = '  +1.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06624'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.36%  '

$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("D22").Value = '28.383.78'
$ws.Range("E22").Value = '  +3.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.473'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.20%  '

$ws.Range("E24").Value = '  +1.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.256'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").Value = '2.087.68'
$ws.Range("E26").Value = '  +1.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.28%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.122'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.77%  '

$ws.Range("E30").Value = '  +3.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '120.37'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9696'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09509'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.598'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.78%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.373'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.47%  '

$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.348'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.93%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06095'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02251'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.359'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.178'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.48%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5947'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.57%  '

$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1875'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.36'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.39%  '

$ws.Range("E45").Value = '  +3.17%  '

$ws.Range("E46").Value = '  +1.99%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.98%  '

$ws.Range("E48").Value = '  +4.90%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06904'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.059'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +15.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '111.66'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.40%  '

